$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 35715670
$ws.Range("I11").Value = 35715670
$ws.Range("K11").Value = 35715670
$ws.Range("M11").Value = -35715530
$ws.Range("H53").Value = 548.8125
$ws.Range("I53").Value = 364.66666
$ws.Range("J53").Value = 785.5714
$ws.Range("K53").Value = 364.66666
$ws.Range("L53").Value = 785.5714
$ws.Range("M53").Value = 272.33334
$ws.Range("N53").Value = -2059.5714
$ws.Range("H55").Value = 1585.7142
$ws.Range("I55").Value = 173.27272
$ws.Range("K55").Value = 173.27272
$ws.Range("M55").Value = 40.72728000000001
$ws.Range("H115").Value = 225
$ws.Range("I115").Value = 225
$ws.Range("K115").Value = 675
$ws.Range("M115").Value = 892
$ws.Range("H135").Value = 14293.111
$ws.Range("I135").Value = 1394.8
$ws.Range("K135").Value = 12553.2
$ws.Range("M135").Value = -10018.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 1000
$ws.Range("J23").Value = 1000
$ws.Range("L23").Value = 1000
$ws.Range("N23").Value = -1518
$ws.Range("H32").Value = 10205481
$ws.Range("I32").Value = 10418090
$ws.Range("K32").Value = 10418090
$ws.Range("M32").Value = -10417803
$ws.Range("H43").Value = 29777.334
$ws.Range("J43").Value = 33333
$ws.Range("L43").Value = 33333
$ws.Range("N43").Value = -33959
$ws.Range("H45").Value = 1962.8667
$ws.Range("I45").Value = 1868.625
$ws.Range("J45").Value = 2070.5715
$ws.Range("K45").Value = 1868.625
$ws.Range("L45").Value = 2070.5715
$ws.Range("M45").Value = -1491.625
$ws.Range("N45").Value = -2824.5715
$ws.Range("H61").Value = 41761340
$ws.Range("I61").Value = 100008550
$ws.Range("K61").Value = 100008550
$ws.Range("M61").Value = -100008338
$ws.Range("H74").Value = 13899546
$ws.Range("I74").Value = 27778708
$ws.Range("K74").Value = 27778708
$ws.Range("M74").Value = -27777834
$ws.Range("H77").Value = 13899546
$ws.Range("I77").Value = 27778708
$ws.Range("K77").Value = 138893540
$ws.Range("M77").Value = -138889172
$ws.Range("H132").Value = 5665.8687
$ws.Range("I132").Value = 3307.0356
$ws.Range("J132").Value = 12270.6
$ws.Range("K132").Value = 9921.106800000001
$ws.Range("L132").Value = 36811.8
$ws.Range("M132").Value = -7391.106800000001
$ws.Range("N132").Value = -41871.8
$ws.Range("H133").Value = 63999.5
$ws.Range("J133").Value = 63999.5
$ws.Range("L133").Value = 63999.5
$ws.Range("N133").Value = -69059.5
$ws.Range("H136").Value = 41761340
$ws.Range("I136").Value = 100008550
$ws.Range("K136").Value = 300025650
$ws.Range("M136").Value = -300023100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 221000
$ws.Range("J70").Value = 221000
$ws.Range("L70").Value = 221000
$ws.Range("N70").Value = -221586
$ws.Range("H73").Value = 221000
$ws.Range("J73").Value = 221000
$ws.Range("L73").Value = 221000
$ws.Range("N73").Value = -223028
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("H134").Value = 67874.75
$ws.Range("I134").Value = 2848.9
$ws.Range("K134").Value = 8546.700000000001
$ws.Range("M134").Value = -6011.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 940648.75
$ws.Range("I31").Value = 14291.417
$ws.Range("K31").Value = 14291.417
$ws.Range("M31").Value = -13996.417
$ws.Range("H34").Value = 940648.75
$ws.Range("I34").Value = 14291.417
$ws.Range("K34").Value = 14291.417
$ws.Range("M34").Value = -14089.417
$ws.Range("H119").Value = 80037.664
$ws.Range("J119").Value = 80037.664
$ws.Range("L119").Value = 80037.664
$ws.Range("N119").Value = -89713.664
$ws.Range("H132").Value = 2480.9678
$ws.Range("I132").Value = 2237.2083
$ws.Range("J132").Value = 3316.7144
$ws.Range("K132").Value = 6711.624899999999
$ws.Range("L132").Value = 9950.143199999999
$ws.Range("M132").Value = -4181.624899999999
$ws.Range("N132").Value = -15010.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 50
$ws.Range("I7").Value = 50
$ws.Range("K7").Value = 150
$ws.Range("M7").Value = -38
$ws.Range("H14").Value = 699.5
$ws.Range("I14").Value = 699.5
$ws.Range("K14").Value = 2098.5
$ws.Range("M14").Value = -1925.5
$ws.Range("H34").Value = 787
$ws.Range("I34").Value = 787
$ws.Range("K34").Value = 2361
$ws.Range("M34").Value = -2277
$ws.Range("H97").Value = 1480.1666
$ws.Range("I97").Value = 992.6
$ws.Range("K97").Value = 2977.8
$ws.Range("M97").Value = -2481.8
$ws.Range("H124").Value = 1553.75
$ws.Range("J124").Value = 700
$ws.Range("L124").Value = 2100
$ws.Range("N124").Value = -11920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4709.5557
$ws.Range("I70").Value = 4709.5557
$ws.Range("K70").Value = 4709.5557
$ws.Range("M70").Value = -4439.5557
$ws.Range("H73").Value = 4709.5557
$ws.Range("I73").Value = 4709.5557
$ws.Range("K73").Value = 4709.5557
$ws.Range("M73").Value = -3773.5557
$ws.Range("H109").Value = 46370.5
$ws.Range("J109").Value = 46370.5
$ws.Range("L109").Value = 46370.5
$ws.Range("N109").Value = -48450.5
$ws.Range("H132").Value = 45457460
$ws.Range("I132").Value = 66669696
$ws.Range("K132").Value = 200009088
$ws.Range("M132").Value = -200006558

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49661.137
$ws.Range("I7").Value = 3472.1177
$ws.Range("K7").Value = 3472.1177
$ws.Range("M7").Value = -3360.1177
$ws.Range("H16").Value = 1704.6923
$ws.Range("I16").Value = 1518.7142
$ws.Range("J16").Value = 1921.6666
$ws.Range("K16").Value = 1518.7142
$ws.Range("L16").Value = 1921.6666
$ws.Range("M16").Value = -1348.7142
$ws.Range("N16").Value = -2261.6666
$ws.Range("H82").Value = 1631
$ws.Range("J82").Value = 1334
$ws.Range("L82").Value = 1334
$ws.Range("N82").Value = -2056
$ws.Range("H85").Value = 1631
$ws.Range("J85").Value = 1334
$ws.Range("L85").Value = 1334
$ws.Range("N85").Value = -3830
$ws.Range("H93").Value = 166670000
$ws.Range("I93").Value = 250002990
$ws.Range("K93").Value = 250002990
$ws.Range("M93").Value = -250001742
$ws.Range("H126").Value = 49661.137
$ws.Range("I126").Value = 3472.1177
$ws.Range("K126").Value = 10416.3531
$ws.Range("M126").Value = -7946.3531
$ws.Range("H127").Value = 155950
$ws.Range("J127").Value = 155950
$ws.Range("L127").Value = 155950
$ws.Range("N127").Value = -165870
$ws.Range("H128").Value = 97473.39999999999
$ws.Range("J128").Value = 97473.39999999999
$ws.Range("L128").Value = 97473.39999999999
$ws.Range("N128").Value = -107433.4
$ws.Range("H129").Value = 78000
$ws.Range("J129").Value = 78000
$ws.Range("L129").Value = 78000
$ws.Range("N129").Value = -88000
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H132").Value = 460053.6
$ws.Range("I132").Value = 436577.44
$ws.Range("K132").Value = 1309732.32
$ws.Range("M132").Value = -1307202.32
$ws.Range("H136").Value = 371501.66
$ws.Range("I136").Value = 7250
$ws.Range("K136").Value = 21750
$ws.Range("M136").Value = -19200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30495
$ws.Range("J40").Value = 30495
$ws.Range("L40").Value = 30495
$ws.Range("N40").Value = -30793
$ws.Range("I62").Value = 6858
$ws.Range("J62").Value = 40004500
$ws.Range("K62").Value = 6858
$ws.Range("L62").Value = 40004500
$ws.Range("M62").Value = -6234
$ws.Range("N62").Value = -40005748
$ws.Range("I65").Value = 6858
$ws.Range("J65").Value = 40004500
$ws.Range("K65").Value = 34290
$ws.Range("L65").Value = 200022500
$ws.Range("M65").Value = -31170
$ws.Range("N65").Value = -200028740
$ws.Range("H93").Value = 81941.336
$ws.Range("J93").Value = 77912
$ws.Range("L93").Value = 77912
$ws.Range("N93").Value = -82904
$ws.Range("H122").Value = 4253.731
$ws.Range("I122").Value = 2083.7778
$ws.Range("K122").Value = 6251.3334
$ws.Range("M122").Value = -3801.3334
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 1833.579
$ws.Range("I132").Value = 1833.579
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5500.737
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2970.737
$ws.Range("N132").Value = $null
$ws.Range("H136").Value = 1313.5454
$ws.Range("I136").Value = 1397.4
$ws.Range("J136").Value = 475
$ws.Range("K136").Value = 4192.200000000001
$ws.Range("L136").Value = 1425
$ws.Range("M136").Value = -1642.200000000001
$ws.Range("N136").Value = -6525
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
